$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.627.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.597.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.597.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.626.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.654"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.291.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.735.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.893"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0504"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
